$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "Realizar una pagina principal en la cual se encuentren los siguientes
# elementos." becomes two runs: the sentence (minus the final period) gets
# struck-through and wrapped in a "_GoBack" bookmark, the trailing "."
# stays a normal run.
$r1 = $d.Content
$r1.Find.Execute(
    "Realizar una pagina principal en la cual se encuentren los siguientes elementos",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Font.StrikeThrough = 1
$d.Bookmarks.Add("_GoBack", $r1)

# --- Change 2 -------------------------------------------------------------
# The old "_GoBack" bookmark that split "El lista" / "do de estos..." is
# removed, and the two strike-through runs are merged back into a single
# run reading "El listado de estos se podrá observar desde la pantalla de
# administración. ".
$r2 = $d.Content
$r2.Find.Execute(
    "El listado de estos se podrá observar desde la pantalla de administración. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El listado de estos se podrá observar desde la pantalla de administración. ", 2)
